$d = $word.ActiveDocument

# Fix 1: "As at {{ doc_generated_date }}" was previously split across three
# runs as "As at {{ " + "doc_generated" + "_date }}". Re-running Find/Replace
# over the same (already correct) rendered text collapses it into one run,
# matching the target formatting.
$d.Content.Find.Execute("As at {{ doc_generated_date }}", $true, $false, $false, $false, $false, $true, 1, $false, "As at {{ doc_generated_date }}", 2)

# Fix 2: the "Authorised Users..." line incorrectly referenced issue_date;
# it should reference doc_generated_date instead. Use the full unique phrase
# so the other, unrelated "Date: {{ issue_date }}" occurrence is untouched.
$d.Content.Find.Execute("Authorised Users of {{ mooring_name }} as at {{ issue_date }}", $true, $false, $false, $false, $false, $true, 1, $false, "Authorised Users of {{ mooring_name }} as at {{ doc_generated_date }}", 2)
